$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Insert a new blank row at row 6 (pushes old rows 6-9 down to 7-10),
#    then sync the table's range to include it (and the 4 rows we will
#    append below) by growing via ListRows.Add().
$ws.Rows.Item(6).Insert() | Out-Null
$lo.ListRows.Add() | Out-Null   # table range -> A1:C10
$lo.ListRows.Add() | Out-Null   # table range -> A1:C11
$lo.ListRows.Add() | Out-Null   # table range -> A1:C12
$lo.ListRows.Add() | Out-Null   # table range -> A1:C13
$lo.ListRows.Add() | Out-Null   # table range -> A1:C14

# 2) Fill the new row 6 ("database" / "facoep " / True) - first new strings.
$ws.Range("A6").Value = "database"
$ws.Range("B6").Value = "facoep "
$ws.Range("C2").Copy($ws.Range("C6"))

# 3) Fill new row 11 ("database" / "Facoep" / False) before touching row 10's
#    Valor cell, so the shared-string append order matches the source edit.
$ws.Range("A11").Value = "database"
$ws.Range("B11").Value = "Facoep"
$ws.Range("C7").Copy($ws.Range("C11"))

# 4) Now update row 10's Valor (old "...Cobranzas/Versión 7" path) to the new
#    MonitoreoCRGs path (A10/C10 already hold the right text after the shift).
$ws.Range("B10").Value = "E:/Personales/Sistemas/Agustin/Reportes BI/2021/MonitoreoCRGs"

# 5) Fill new row 12 (host / 10.22.0.142 / False) with the Arial 10pt font.
$ws.Range("A12").Value = "host"
$ws.Range("B12").Value = "10.22.0.142"
$ws.Range("C7").Copy($ws.Range("C12"))
$ws.Range("B12").Font.Name = "Arial"
$ws.Range("B12").Font.Size = 10

# 6) Fill new row 13 (user / postgres / False), same font via format copy.
$ws.Range("A13").Value = "user"
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("C7").Copy($ws.Range("C13"))
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

# 7) Fill new row 14 (password / serveradmin / False), same font.
$ws.Range("A14").Value = "password"
$ws.Range("B14").Value = "serveradmin"
$ws.Range("C7").Copy($ws.Range("C14"))
$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# 8) Page setup orientation (best-effort match of the worksheet pageSetup
#    element added in the diff).
$ws.PageSetup.Orientation = 1

# 9) Selection bookkeeping to match the saved cursor position in the diff.
$ws.Range("B22").Select() | Out-Null
